$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4
$ws.Range("C4").Value = 1893.4
$ws.Range("D4").Value = 1885.09320921534
$ws.Range("E4").Value = 1907.3318350055
$ws.Range("F4").Value = 1930.4539814257
$ws.Range("G4").Value = 1954.158923269
$ws.Range("H4").Value = 1968.11289188136
$ws.Range("I4").Value = 1980.71518106142
$ws.Range("J4").Value = 1995.67745272291
$ws.Range("K4").Value = 2009.40022998616

# Row 5
$ws.Range("C5").Value = 3139.2
$ws.Range("D5").Value = 3167.88900514805
$ws.Range("E5").Value = 3195.56015057414
$ws.Range("F5").Value = 3223.36403051996
$ws.Range("G5").Value = 3250.16069012155
$ws.Range("H5").Value = 3276.93438119768
$ws.Range("I5").Value = 3304.82349595586
$ws.Range("J5").Value = 3332.32958624848
$ws.Range("K5").Value = 3360.14586941614

# Row 10
$ws.Range("D10").Value = 647.018687447845
$ws.Range("E10").Value = 647.926324245673
$ws.Range("F10").Value = 650.230137472713
$ws.Range("G10").Value = 652.537182485326
$ws.Range("H10").Value = 656.890013221912
$ws.Range("I10").Value = 661.262575183405
$ws.Range("J10").Value = 665.654957810736
$ws.Range("K10").Value = 670.06725095027

# Row 11
$ws.Range("C11").Value = 919.4
$ws.Range("D11").Value = 927.689732874973
$ws.Range("E11").Value = 928.991094574193
$ws.Range("F11").Value = 932.294281821561
$ws.Range("G11").Value = 935.602102774806
$ws.Range("H11").Value = 941.843153399173
$ws.Range("I11").Value = 948.11249448117
$ws.Range("J11").Value = 954.410254260423
$ws.Range("K11").Value = 960.736561557861

# Row 17
$ws.Range("C17").Value = 2152.299
$ws.Range("D17").Value = 2160.527
$ws.Range("E17").Value = 2202.093889974
$ws.Range("F17").Value = 2208.793889974
$ws.Range("G17").Value = 2215.493889974
$ws.Range("H17").Value = 2205.954889974
$ws.Range("I17").Value = 2249.7275234174
$ws.Range("J17").Value = 2255.6995234174
$ws.Range("K17").Value = 2263.3995234174

# Row 19
$ws.Range("C19").Value = 4465.4
$ws.Range("D19").Value = 4505.15020255048
$ws.Range("E19").Value = 4578.04931402487
$ws.Range("F19").Value = 4609.71515475437
$ws.Range("G19").Value = 4641.62782503496
$ws.Range("H19").Value = 4673.78951251855
$ws.Range("I19").Value = 4755.93289067343
$ws.Range("J19").Value = 4788.84607274818
$ws.Range("K19").Value = 4822.01619301344

# Row 20
$ws.Range("C20").Value = 2483.3
$ws.Range("D20").Value = 2522.33520327819
$ws.Range("E20").Value = 2548.44772380158
$ws.Range("F20").Value = 2574.48894008494
$ws.Range("G20").Value = 2599.52236547798
$ws.Range("H20").Value = 2624.24639104833
$ws.Range("I20").Value = 2646.38928336033
$ws.Range("J20").Value = 2669.39124256
$ws.Range("K20").Value = 2692.90294087333

# Row 21
$ws.Range("C21").Value = 483.1
$ws.Range("D21").Value = 492.801912246684
$ws.Range("E21").Value = 502.698664280664
$ws.Range("F21").Value = 512.794168994753
$ws.Range("G21").Value = 523.092417862893
$ws.Range("H21").Value = 521.771537588453
$ws.Range("I21").Value = 520.453992718313
$ws.Range("J21").Value = 519.139774830117
$ws.Range("K21").Value = 517.828875522779

# Row 22
$ws.Range("C22").Value = 164.8
$ws.Range("D22").Value = 163.15466849793
$ws.Range("E22").Value = 164.648939314487
$ws.Range("F22").Value = 166.13773661157
$ws.Range("G22").Value = 167.08465548067
$ws.Range("H22").Value = 167.166758272789
$ws.Range("I22").Value = 166.433306663197
$ws.Range("J22").Value = 166.312889234756
$ws.Range("K22").Value = 166.756244312196
